# daily auto push: 2026-02-01 18:52 UTC
# Insert two new daily-sample rows (2026/02/01 22:00 and 2026/02/02 02:00)
# right before the existing "2026/12/29" block, shifting all subsequent
# rows down by two and extending the used range to A1:D807.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every row from 764 downward two rows lower, then populate the
# two freshly inserted rows (new 764 & 765).
$ws.Rows.Item(764).Resize(2).Insert()

# Force the date-like text into the cells as literal strings (otherwise
# Excel's smart entry would reinterpret "2026/02/01" as a date serial),
# then drop the temporary Text format so the cells end up styled the same
# as the rest of the data rows (no explicit style index).
$ws.Range("A764:A765").NumberFormat = "@"

$ws.Range("A764").Value = "2026/02/01"
$ws.Range("B764").Value = "日"
$ws.Range("C764").Value = 22
$ws.Range("D764").Value = 152

$ws.Range("A765").Value = "2026/02/02"
$ws.Range("B765").Value = "月"
$ws.Range("C765").Value = 2
$ws.Range("D765").Value = 160

$ws.Range("A764:D765").ClearFormats()
